$wb = $excel.ActiveWorkbook

# Rename the first sheet from "sample1" to "sheet1"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "sheet1"

# Update header cells: A1 "Month" -> "sh", B1 "Average" -> " Average"
$ws1.Range("A1").Value = "sh"
$ws1.Range("B1").Value = " Average"

# Add a new worksheet named "ask" after sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ask"

# Keep sheet1 as the active/selected sheet, with D27 selected
$ws1.Activate() | Out-Null
$ws1.Range("D27").Select() | Out-Null
